$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header fields: account holder name ---
$ws.Range("C2").Value = "Hartmut"

# Card number is a 16-digit value that must stay exact text (Excel's
# General number format would otherwise round it to 15 significant
# digits). Briefly force Text format for the assignment, then restore
# the original "General" format so the cell keeps its original style.
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "2570314725427075"
$ws.Range("B3").NumberFormat = "general"

$ws.Range("C3").Value = "Mohaupt"

# --- Opening balance line ---
$ws.Range("D5").Value = "KONTOSTAND AM 03.01.2025"

# --- Transaction row 6 (re-dated + new description/amount) ---
$ws.Range("B6").Value = "07.01."
$ws.Range("C6").Value = "08.01."
$ws.Range("D6").Value = "KARTENZAHLUNG SHELL TANKSTELLE"
$ws.Range("E6").Value = "83,42-"

# --- Transaction row 7 ---
$ws.Range("B7").Value = "09.01."
$ws.Range("C7").Value = "10.01."
$ws.Range("D7").Value = "BEITRAG Allianz SE K-41896920"
$ws.Range("E7").Value = "53,48-"

# --- Transaction row 8 ---
$ws.Range("B8").Value = "13.01."
$ws.Range("C8").Value = "14.01."
$ws.Range("D8").Value = "ABSCHLAG STROM Stadtwerke Rosenheim 63257966"
$ws.Range("E8").Value = "84,36-"

# --- Rows 9 & 10 previously held transactions; they are now blank spacer
#     rows (matching the style already used by the blank row 11) ---
$ws.Range("B9:D10").ClearContents()
$ws.Range("E9").Value = ""
$ws.Range("E10").Value = ""

# Re-apply the blank-row formatting used elsewhere (E11 already carries
# the target look for E10; E10's look - after flipping to centered - is
# what E9 needs), so no new styles need to be fabricated.
$ws.Range("E11").Copy()
$ws.Range("E10").PasteSpecial(-4122)

$ws.Range("E10").Copy()
$ws.Range("E9").PasteSpecial(-4122)
$ws.Range("E9").HorizontalAlignment = -4108

# --- Closing balance line ---
$ws.Range("D12").Value = "KONTOSTAND AM 15.01.2025"
$ws.Range("E12").Value = "221,26-"

# --- Next statement date ---
$ws.Range("C13").Value = "IHR NAECHSTER ABRECHNUNGSTERMIN 23.01.2025"
